$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to Text before writing numeric-looking strings,
# so values like "0.9999" or "242.09" are stored as text (matching the
# original inline-string cells) instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.328.32'
$ws.Range("E2").Value = '  +0.29%  '

$ws.Range("D3").Value = '1.877.25'
$ws.Range("E3").Value = '  +0.43%  '

$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '0.7121'
$ws.Range("E5").Value = '  +0.26%  '

$ws.Range("D6").Value = '242.09'
$ws.Range("E6").Value = '  +0.29%  '

$ws.Range("D7").Value = '0.9999'
$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").Value = '0.07881'
$ws.Range("E8").Value = '  +2.58%  '

$ws.Range("E9").Value = '  +0.61%  '

$ws.Range("D10").Value = '25.32'
$ws.Range("E10").Value = '  +1.52%  '

$ws.Range("D11").Value = '0.08391'
$ws.Range("E11").Value = '  +0.24%  '

$ws.Range("D12").Value = '1.867.04'
$ws.Range("E12").Value = '  -0.45%  '

$ws.Range("D13").Value = '5.254'
$ws.Range("E13").Value = '  +1.00%  '

$ws.Range("D14").Value = '0.7193'
$ws.Range("E14").Value = '  +1.46%  '

$ws.Range("D15").Value = '91.39'
$ws.Range("E15").Value = '  +0.36%  '

$ws.Range("D16").Value = '6.229'
$ws.Range("E16").Value = '  +5.12%  '

$ws.Range("D17").Value = '0.000008353'
$ws.Range("E17").Value = '  +1.19%  '

$ws.Range("D18").Value = '29.329.40'

$ws.Range("D19").Value = '240.97'
$ws.Range("E19").Value = '  -0.52%  '

$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '13.25'
$ws.Range("E20").Value = '  +0.80%  '

$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.120.75'
$ws.Range("E21").Value = '  -0.31%  '

$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  +0.04%  '

$ws.Range("D23").Value = '7.796'
$ws.Range("E23").Value = '  -0.11%  '

$ws.Range("E24").Value = '  -0.02%  '

$ws.Range("E25").Value = '  -1.66%  '

$ws.Range("D26").Value = '162.91'
$ws.Range("E26").Value = '  -0.05%  '

$ws.Range("D27").Value = '9.064'
$ws.Range("E27").Value = '  +0.80%  '

$ws.Range("E28").Value = '  +0.29%  '

$ws.Range("D29").Value = '1.510'
$ws.Range("E29").Value = '  +0.41%  '

$ws.Range("D30").Value = '4.424'
$ws.Range("E30").Value = '  +0.54%  '

$ws.Range("D31").Value = '4.344'

$ws.Range("D32").Value = '1.224'
$ws.Range("E32").Value = '  -4.34%  '

$ws.Range("D33").Value = '0.05364'
$ws.Range("E33").Value = '  +2.40%  '

$ws.Range("D34").Value = '1.951'
$ws.Range("E34").Value = '  +1.61%  '

$ws.Range("D35").Value = '1.181'
$ws.Range("E35").Value = '  +1.02%  '

$ws.Range("D36").Value = '0.7465'
$ws.Range("E36").Value = '  -0.51%  '

$ws.Range("D37").Value = '2.686'
$ws.Range("E37").Value = '  +0.07%  '

$ws.Range("D38").Value = '1.306.95'
$ws.Range("E38").Value = '  +13.25%  '

$ws.Range("E39").Value = '  +1.55%  '

$ws.Range("D40").Value = '2.741'
$ws.Range("E40").Value = '  +0.98%  '

$ws.Range("D41").Value = '6.528'
$ws.Range("E41").Value = '  +2.76%  '

$ws.Range("D42").Value = '110.74'
$ws.Range("E42").Value = '  +6.22%  '

$ws.Range("D43").Value = '0.8934'
$ws.Range("E43").Value = '  +0.99%  '

$ws.Range("D44").Value = '73.14'
$ws.Range("E44").Value = '  +0.23%  '

$ws.Range("D45").Value = '0.00000000133'
$ws.Range("E45").Value = '  +14.55%  '

$ws.Range("D46").Value = '0.9998'
$ws.Range("E46").Value = '  -0.06%  '

$ws.Range("D47").Value = '2.029.46'
$ws.Range("E47").Value = '  +0.16%  '

$ws.Range("E48").Value = '  +0.81%  '

$ws.Range("D49").Value = '0.5189'
$ws.Range("E49").Value = '  +0.11%  '

$ws.Range("D50").Value = '9.486'
$ws.Range("E50").Value = '  +1.35%  '

$ws.Range("D51").Value = '0.4363'
$ws.Range("E51").Value = '  +1.82%  '

# Restore the column to the workbook's default (Normal) style/number format
# so the saved cells carry no explicit style index, same as the source file.
$ws.Range("D2:D51").Style = "Normal"
